$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trial0")

# Replace the calculated-column formulas in F5:F20 with plain (hard-coded)
# numeric values, matching the author's "forgot to save" correction where
# the live formula results were overwritten by manually typed numbers.
$values = @{
    5  = 9.8000000000000007
    6  = 19
    7  = 29
    8  = 36
    9  = 38
    10 = 48
    11 = 58
    12 = 56
    13 = 78
    14 = 79
    15 = 90
    16 = 110
    17 = 130
    18 = 140
    19 = 145
    20 = 155
}

foreach ($row in $values.Keys) {
    $cell = $ws.Range("F$row")
    $cell.Value = $values[$row]
}

$wb.Save()
